$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark exercises 5 through 10 (rows 6-11, column B) as done (TRUE)
$ws.Range("B6:B11").Value = $true

# Underline the statement of exercise 7 (row 8, column A) to highlight it
$ws.Range("A8").Font.Underline = $true

# Update the view: scroll so row 2 is the top row, and select A12
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A12").Select()

# Configure page setup for printing (A4, portrait)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
